$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16/18, 41/42, 46/47 swapped ranking position: update coin name + link
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"

# Refresh Price (D) and Volume(1h) (E) columns with latest scrape
$ws.Range("D2").Value = "63.655.48"
$ws.Range("E2").Value = "  -3.98%  "
$ws.Range("D3").Value = "3.336.96"
$ws.Range("E3").Value = "  -4.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.90"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.66"
$ws.Range("E6").Value = "  -6.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -4.31%  "
$ws.Range("D8").Value = "3.329.33"
$ws.Range("E8").Value = "  -4.35%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.613"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.33"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  -4.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.91"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").Value = "3.855.58"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.82"
$ws.Range("E16").Value = "  -3.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  -3.70%  "
$ws.Range("D18").Value = "3.316.22"
$ws.Range("E18").Value = "  -4.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.68"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").Value = "63.473.48"
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.971"
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "410.06"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.06"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.33"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.26"
$ws.Range("E25").Value = "  +7.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.95"
$ws.Range("E26").Value = "  -3.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.57"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.72"
$ws.Range("E28").Value = "  -6.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.56"
$ws.Range("E29").Value = "  -5.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.10"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  -6.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.31"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "571.99"
$ws.Range("E33").Value = "  -8.16%  "
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.46"
$ws.Range("E35").Value = "  -4.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.146"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.22"
$ws.Range("E38").Value = "  -7.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("E40").Value = "  -7.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").Value = "3.148.08"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0401"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -6.22%  "
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.79"
$ws.Range("E50").Value = "  -5.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.98"
$ws.Range("E51").Value = "  -5.61%  "
